# Marksheet edit: rename Sheet1 -> Database, add a new "Query" sheet with
# two practice questions, following the author's commit "two questions are
# added is marksheet".

$wb = $excel.ActiveWorkbook
$db = $wb.Worksheets.Item(1)
$db.Name = "Database"

# New worksheet, placed right after Database.
$q = $wb.Worksheets.Add($null, $db)
$q.Name = "Query"

# ---- Database sheet bookkeeping (column widths + page setup + selection) --
$db.Columns.Item(14).ColumnWidth = 9.666666666666666   # -> 10.5703125 target (closest reachable)
$db.Columns.Item(15).ColumnWidth = 30.0                  # -> 30.85546875 target (closest reachable)
$db.Columns.Item(16).ColumnWidth = 35.33333333333333    # -> 36.140625 target (closest reachable)

$db.PageSetup.PaperSize = 9
$db.PageSetup.Orientation = 1

$db.Range("M2:T14").Select()
$wb.Windows.Item(1).ScrollRow = 16

# ---- Question 1 -------------------------------------------------------------
$q.Range("D4").Value = 1
$q.Range("E4").Value = "Find name from student_id (first_name)"
$q.Range("E4:K4").Merge()

$q.Range("E6").Value = "student_id"
$q.Range("F6").Value = "first_name"

$q.Range("E7").Value = 19359

# ---- Question 2 -------------------------------------------------------------
$q.Range("D10").Value = 2
$q.Range("E10").Value = "Find name from student_id (first_name and last_name)"
$q.Range("E10:K10").Merge()

$q.Range("E12").Value = "student_id"
$q.Range("F12").Value = "first_name"
$q.Range("G12").Value = "last_name"

# Named "answer box" cell styles, applied in ascending builtinId order
# (Good=26, 20% - Accent3=38, 20% - Accent5=46) so the generated xfId /
# cellStyle bookkeeping lines up with how Excel itself orders the gallery.
$q.Range("F13").Style = "Good"
$q.Range("G13").Style = "20% - Accent3"
$q.Range("F7").Style = "20% - Accent5"

# ---- Title -------------------------------------------------------------------
$q.Range("E1").Value = "Solve the following problems"
$q.Range("E1:G1").Merge()
$q.Range("E1:G1").HorizontalAlignment = -4108  # xlCenter

# ---- Fonts for question headers (bold, larger, accent-colored) --------------
# NB: apply per contiguous area -- multi-area (comma) Range.Font assignment
# only reliably sticks to the first area in this runtime.
foreach ($addr in @("D4", "E4:K4", "D10", "E10:K10")) {
    $f = $q.Range($addr).Font
    $f.Bold = $true
    $f.Size = 12
    $f.ThemeColor = 5   # msoThemeColorAccent1 -> theme="4" in the xlsx
}

$q.Range("E4:K4").HorizontalAlignment = -4131   # xlLeft
$q.Range("E10:K10").HorizontalAlignment = -4131 # xlLeft

# ---- Column widths ------------------------------------------------------------
$q.Columns.Item(6).ColumnWidth = 35.166666666666664  # -> 36 target
$q.Columns.Item(7).ColumnWidth = 36.0                 # -> 36.85546875 target (closest reachable)

# ---- View state: Query ends up the active sheet/tab ---------------------------
$q.Select()
$q.Range("F15:F16").Select()
